$wb = $excel.ActiveWorkbook

# --- Sheet: LP1912 ---
$ws = $wb.Worksheets.Item("LP1912")
$ws.Range("A110:E115").EntireRow.Insert()
$ws.Range("A2").Value = "Última actualización: 08:56:29"
$ws.Range("A3").Value = "Total filas: 122"
$ws.Cells.Item(110,1).Value = "08:56:29"
$ws.Cells.Item(110,2).Value = "09:34"
$ws.Cells.Item(110,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(110,4).Value = 38
$ws.Cells.Item(110,5).Value = "LP1912"
$ws.Cells.Item(111,1).Value = "08:56:29"
$ws.Cells.Item(111,2).Value = "09:34"
$ws.Cells.Item(111,3).Value = "16_SANTA ANA"
$ws.Cells.Item(111,4).Value = 38
$ws.Cells.Item(111,5).Value = "LP1912"
$ws.Cells.Item(112,1).Value = "08:32:09"
$ws.Cells.Item(112,2).Value = "09:35"
$ws.Cells.Item(112,3).Value = "16_SANTA ANA"
$ws.Cells.Item(112,4).Value = 63
$ws.Cells.Item(112,5).Value = "LP1912"
$ws.Cells.Item(113,1).Value = "08:48:08"
$ws.Cells.Item(113,2).Value = "09:35"
$ws.Cells.Item(113,3).Value = "23_HERNANDEZ"
$ws.Cells.Item(113,4).Value = 47
$ws.Cells.Item(113,5).Value = "LP1912"
$ws.Cells.Item(114,1).Value = "07:50:16"
$ws.Cells.Item(114,2).Value = "09:42"
$ws.Cells.Item(114,3).Value = "215C_EL PATO"
$ws.Cells.Item(114,4).Value = 112
$ws.Cells.Item(114,5).Value = "LP1912"
$ws.Cells.Item(115,1).Value = "08:02:22"
$ws.Cells.Item(115,2).Value = "09:43"
$ws.Cells.Item(115,3).Value = "14_ABASTO"
$ws.Cells.Item(115,4).Value = 101
$ws.Cells.Item(115,5).Value = "LP1912"
$ws.Cells.Item(116,1).Value = "07:50:16"
$ws.Cells.Item(116,2).Value = "09:44"
$ws.Cells.Item(116,3).Value = "14_ABASTO"
$ws.Cells.Item(116,4).Value = 114
$ws.Cells.Item(116,5).Value = "LP1912"
$ws.Cells.Item(117,1).Value = "08:32:09"
$ws.Cells.Item(117,2).Value = "09:52"
$ws.Cells.Item(117,3).Value = "15_ABASTO"
$ws.Cells.Item(117,4).Value = 80
$ws.Cells.Item(117,5).Value = "LP1912"
$ws.Cells.Item(118,1).Value = "08:56:29"
$ws.Cells.Item(118,2).Value = "09:53"
$ws.Cells.Item(118,3).Value = "10_OLMOS"
$ws.Cells.Item(118,4).Value = 57
$ws.Cells.Item(118,5).Value = "LP1912"
$ws.Cells.Item(119,1).Value = "08:56:29"
$ws.Cells.Item(119,2).Value = "10:10"
$ws.Cells.Item(119,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(119,4).Value = 74
$ws.Cells.Item(119,5).Value = "LP1912"
$ws.Cells.Item(120,1).Value = "08:32:09"
$ws.Cells.Item(120,2).Value = "10:11"
$ws.Cells.Item(120,3).Value = "16_P MOR-SANTA ANA"
$ws.Cells.Item(120,4).Value = 99
$ws.Cells.Item(120,5).Value = "LP1912"
$ws.Cells.Item(121,1).Value = "08:32:09"
$ws.Cells.Item(121,2).Value = "10:21"
$ws.Cells.Item(121,3).Value = "26_HERNANDEZ"
$ws.Cells.Item(121,4).Value = 109
$ws.Cells.Item(121,5).Value = "LP1912"
$ws.Cells.Item(122,1).Value = "08:32:09"
$ws.Cells.Item(122,2).Value = "10:22"
$ws.Cells.Item(122,3).Value = "17_ROMERO"
$ws.Cells.Item(122,4).Value = 110
$ws.Cells.Item(122,5).Value = "LP1912"
$ws.Cells.Item(123,1).Value = "08:56:29"
$ws.Cells.Item(123,2).Value = "10:26"
$ws.Cells.Item(123,3).Value = "215A_EL PATO"
$ws.Cells.Item(123,4).Value = 90
$ws.Cells.Item(123,5).Value = "LP1912"
$ws.Cells.Item(124,1).Value = "08:32:09"
$ws.Cells.Item(124,2).Value = "10:27"
$ws.Cells.Item(124,3).Value = "215A_EL PATO"
$ws.Cells.Item(124,4).Value = 115
$ws.Cells.Item(124,5).Value = "LP1912"
$ws.Cells.Item(125,1).Value = "08:48:08"
$ws.Cells.Item(125,2).Value = "10:42"
$ws.Cells.Item(125,3).Value = "17_ROMERO"
$ws.Cells.Item(125,4).Value = 114
$ws.Cells.Item(125,5).Value = "LP1912"
$ws.Cells.Item(126,1).Value = "08:56:29"
$ws.Cells.Item(126,2).Value = "10:43"
$ws.Cells.Item(126,3).Value = "14_ABASTO"
$ws.Cells.Item(126,4).Value = 107
$ws.Cells.Item(126,5).Value = "LP1912"
$ws.Cells.Item(127,1).Value = "08:48:08"
$ws.Cells.Item(127,2).Value = "10:44"
$ws.Cells.Item(127,3).Value = "14_ABASTO"
$ws.Cells.Item(127,4).Value = 116
$ws.Cells.Item(127,5).Value = "LP1912"

# --- Sheet: LP1912-215 ---
$ws = $wb.Worksheets.Item("LP1912-215")
$ws.Range("A21:E21").EntireRow.Insert()
$ws.Range("A2").Value = "Última actualización: 08:56:29"
$ws.Range("A3").Value = "Total filas: 17"
$ws.Cells.Item(21,1).Value = "08:56:29"
$ws.Cells.Item(21,2).Value = "10:26"
$ws.Cells.Item(21,3).Value = "215A_EL PATO"
$ws.Cells.Item(21,4).Value = 90
$ws.Cells.Item(21,5).Value = "LP1912"
$ws.Cells.Item(22,1).Value = "08:32:09"
$ws.Cells.Item(22,2).Value = "10:27"
$ws.Cells.Item(22,3).Value = "215A_EL PATO"
$ws.Cells.Item(22,4).Value = 115
$ws.Cells.Item(22,5).Value = "LP1912"

# --- Sheet: 6203-6173 ---
$ws = $wb.Worksheets.Item("6203-6173")
$ws.Range("A2").Value = "Última actualización: 08:56:29"
$ws.Range("A3").Value = "Total filas: 25"
$ws.Cells.Item(30,1).Value = "08:56:29"
$ws.Cells.Item(30,2).Value = "10:54"
$ws.Cells.Item(30,3).Value = "215A_LA PLATA"
$ws.Cells.Item(30,4).Value = 118
$ws.Cells.Item(30,5).Value = "L6173"
